# 9.c.1 table: add a "2022" column (column N) to the right of the existing
# 2012-2021 columns (D..M), mirroring the formatting already used by the
# matching cell one column to the left (column M) and filling in the new
# figures for the three data rows (2G/3G/4G coverage).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-ColumnNFormat($row) {
    $ws.Cells.Item($row, 13).Copy() | Out-Null          # column M (source)
    $ws.Cells.Item($row, 14).PasteSpecial($xlPasteFormats) | Out-Null  # column N
    $excel.CutCopyMode = 0
}

# Row 3: blank spacer cell under the thick bottom border - format only.
Copy-ColumnNFormat 3

# Row 4: header row - new year label.
Copy-ColumnNFormat 4
$ws.Cells.Item(4, 14).Value = 2022

# Row 5: 2G coverage for 2022.
Copy-ColumnNFormat 5
$ws.Cells.Item(5, 14).Value = 98.8

# Row 6: 3G coverage for 2022.
Copy-ColumnNFormat 6
$ws.Cells.Item(6, 14).Value = 98

# Row 7: 4G coverage for 2022.
Copy-ColumnNFormat 7
$ws.Cells.Item(7, 14).Value = 96.9

# Leave the saved selection where the workbook shows it - one cell to the
# right of the freshly-typed 2022 header.
$ws.Range("O4").Select()
